$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column group (D:U) gets an explicit width matching the existing pattern
$ws.Range("D1:U15").ColumnWidth = 8.7109375

# Row 3: bottom border style cells extend to T3:U3 (style 4 - thick bottom border)
$ws.Range("T3").Value = $null
$ws.Range("U3").Value = $null

# Row 4: add year headers 2023/2024 (style 6)
$ws.Range("T4").Value = 2023
$ws.Range("U4").Value = 2024

# Row 5: add new data values (style 8) and resize the row height
$ws.Range("T5").Value = 10.8
$ws.Range("U5").Value = 6.5
$ws.Rows.Item(5).RowHeight = 41.25

# Copy styles from existing neighboring cells (S column) across to T/U for rows 3-5
$ws.Range("S3").Copy() | Out-Null
$ws.Range("T3:U3").PasteSpecial(-4122) | Out-Null

$ws.Range("S4").Copy() | Out-Null
$ws.Range("T4:U4").PasteSpecial(-4122) | Out-Null

$ws.Range("S5").Copy() | Out-Null
$ws.Range("T5:U5").PasteSpecial(-4122) | Out-Null

# re-set values since PasteSpecial(formats) shouldn't touch them, but ensure correctness
$ws.Range("T4").Value = 2023
$ws.Range("U4").Value = 2024
$ws.Range("T5").Value = 10.8
$ws.Range("U5").Value = 6.5

# Clear the clipboard marching ants / selection leftover from copy
$excel.CutCopyMode = $false

# Reset worksheet selection to default (A1) to match removal of saved selection state
$ws.Range("A1").Select()
